$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.021714244151236
$ws.Range("D2").Value = 1.025837737684364
$ws.Range("E2").Value = 1.032018076925368
$ws.Range("F2").Value = 1.041893119061836
$ws.Range("I2").Value = 1.028200768471436
$ws.Range("J2").Value = 1.026904155658324
$ws.Range("K2").Value = 1.028662139382646
$ws.Range("L2").Value = 1.034824519679233
$ws.Range("M2").Value = 1.044671346844554
$ws.Range("N2").Value = 1.013042837993227
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.022591839057063
$ws.Range("D3").Value = 1.026458477726111
$ws.Range("E3").Value = 1.032871157811663
$ws.Range("F3").Value = 1.042998839684944
$ws.Range("I3").Value = 1.028318061537426
$ws.Range("J3").Value = 1.027419862389473
$ws.Range("K3").Value = 1.029090773356266
$ws.Range("L3").Value = 1.0354861664013
$ws.Range("M3").Value = 1.045587014308245
$ws.Range("N3").Value = 1.013213480071132
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.023159819975978
$ws.Range("D4").Value = 1.02685967839922
$ws.Range("E4").Value = 1.033423730624165
$ws.Range("F4").Value = 1.04371539572731
$ws.Range("I4").Value = 1.028391943023667
$ws.Range("J4").Value = 1.027753059549461
$ws.Range("K4").Value = 1.029366979940036
$ws.Range("L4").Value = 1.035914227949615
$ws.Range("M4").Value = 1.046180004296005
$ws.Range("N4").Value = 1.013323711471624
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.023398626336766
$ws.Range("D5").Value = 1.027028231471929
$ws.Range("E5").Value = 1.03365616800296
$ws.Range("F5").Value = 1.044016893640168
$ws.Range("I5").Value = 1.028422519552879
$ws.Range("J5").Value = 1.027893015003046
$ws.Range("K5").Value = 1.02948282104907
$ws.Range("L5").Value = 1.036094167679245
$ws.Range("M5").Value = 1.046429414038582
$ws.Range("N5").Value = 1.013370007990928
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.023438724536165
$ws.Range("D6").Value = 1.027056525676727
$ws.Range("E6").Value = 1.033695203186073
$ws.Range("F6").Value = 1.04406753155186
$ws.Range("I6").Value = 1.028427625123673
$ws.Range("J6").Value = 1.027916507035489
$ws.Range("K6").Value = 1.02950225505178
$ws.Range("L6").Value = 1.036124379304306
$ws.Range("M6").Value = 1.046471297848759
$ws.Range("N6").Value = 1.013377778744665
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.023163010812894
$ws.Range("D7").Value = 1.026861931054591
$ws.Range("E7").Value = 1.033426835931925
$ws.Range("F7").Value = 1.043719423345649
$ws.Range("I7").Value = 1.028392353488674
$ws.Range("J7").Value = 1.027754930116392
$ws.Range("K7").Value = 1.029368528901222
$ws.Range("L7").Value = 1.035916632382552
$ws.Range("M7").Value = 1.04618333646371
$ws.Range("N7").Value = 1.013324330264329
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.022010806268796
$ws.Range("D8").Value = 1.026047613869989
$ws.Range("E8").Value = 1.03230626031363
$ws.Range("F8").Value = 1.042266578117554
$ws.Range("I8").Value = 1.028240824910622
$ws.Range("J8").Value = 1.027078543984761
$ws.Range("K8").Value = 1.028807235250598
$ws.Range("L8").Value = 1.0350481397701
$ws.Range("M8").Value = 1.044980698489189
$ws.Range("N8").Value = 1.013100545457042
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.019981431230811
$ws.Range("D9").Value = 1.024609235065134
$ws.Range("E9").Value = 1.030336103050927
$ws.Range("F9").Value = 1.039714801570973
$ws.Range("I9").Value = 1.027958417505298
$ws.Range("J9").Value = 1.025882884284767
$ws.Range("K9").Value = 1.027809424282063
$ws.Range("L9").Value = 1.033517268568137
$ws.Range("M9").Value = 1.042865308618442
$ws.Range("N9").Value = 1.012704803976593
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.018629223940116
$ws.Range("D10").Value = 1.023648094344102
$ws.Range("E10").Value = 1.029025724452769
$ws.Range("F10").Value = 1.038019270560534
$ws.Range("I10").Value = 1.027759845307528
$ws.Range("J10").Value = 1.025083295781817
$ws.Range("K10").Value = 1.027138413199106
$ws.Range("L10").Value = 1.032496427768453
$ws.Range("M10").Value = 1.041457669999093
$ws.Range("N10").Value = 1.012440053769859
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.018043884986323
$ws.Range("D11").Value = 1.023231402966648
$ws.Range("E11").Value = 1.02845905703156
$ws.Range("F11").Value = 1.037286439437866
$ws.Range("I11").Value = 1.027671428477135
$ws.Range("J11").Value = 1.02473648877481
$ws.Range("K11").Value = 1.026846496367999
$ws.Range("L11").Value = 1.032054343102381
$ws.Range("M11").Value = 1.0408487804905
$ws.Range("N11").Value = 1.012325199615195
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.017826491464394
$ws.Range("D12").Value = 1.023076550033601
$ws.Range("E12").Value = 1.028248683163292
$ws.Range("F12").Value = 1.03701443627499
$ws.Range("I12").Value = 1.027638221669633
$ws.Range("J12").Value = 1.024607583041317
$ws.Range("K12").Value = 1.026737861604441
$ws.Range("L12").Value = 1.031890126032672
$ws.Range("M12").Value = 1.040622706957577
$ws.Range("N12").Value = 1.012282505625365
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.01787312185518
$ws.Range("D13").Value = 1.023109769891118
$ws.Range("E13").Value = 1.028293804000962
$ws.Range("F13").Value = 1.037072772678994
$ws.Range("I13").Value = 1.027645361140715
$ws.Range("J13").Value = 1.024635237642794
$ws.Range("K13").Value = 1.026761173336291
$ws.Range("L13").Value = 1.031925351459466
$ws.Range("M13").Value = 1.040671196178079
$ws.Range("N13").Value = 1.012291665075618
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.018025914604534
$ws.Range("D14").Value = 1.02321860429987
$ws.Range("E14").Value = 1.02844166517761
$ws.Range("F14").Value = 1.03726395142648
$ws.Range("I14").Value = 1.027668691030585
$ws.Range("J14").Value = 1.02472583513899
$ws.Range("K14").Value = 1.026837520736444
$ws.Range("L14").Value = 1.032040769008351
$ws.Range("M14").Value = 1.040830091235706
$ws.Range("N14").Value = 1.012321671168361
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.018120058925332
$ws.Range("D15").Value = 1.023285650853295
$ws.Range("E15").Value = 1.028532782157986
$ws.Range("F15").Value = 1.037381769880936
$ws.Range("I15").Value = 1.027683017021596
$ws.Range("J15").Value = 1.024781643856691
$ws.Range("K15").Value = 1.026884533911313
$ws.Range("L15").Value = 1.032111880651652
$ws.Range("M15").Value = 1.04092800435887
$ws.Range("N15").Value = 1.012340154676683
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.018668074748758
$ws.Range("D16").Value = 1.023675738145677
$ws.Range("E16").Value = 1.029063347894142
$ws.Range("F16").Value = 1.038067934560786
$ws.Range("I16").Value = 1.027765662032011
$ws.Range("J16").Value = 1.025106300084319
$ws.Range("K16").Value = 1.027157758106464
$ws.Range("L16").Value = 1.032525766418207
$ws.Range("M16").Value = 1.041498093283502
$ws.Range("N16").Value = 1.01244767174945
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.019011878368959
$ws.Range("D17").Value = 1.023920293883821
$ws.Range("E17").Value = 1.029396355429605
$ws.Range("F17").Value = 1.038498708246974
$ws.Range("I17").Value = 1.027816852033492
$ws.Range("J17").Value = 1.02530979377639
$ws.Range("K17").Value = 1.027328779904833
$ws.Range("L17").Value = 1.032785372276715
$ws.Range("M17").Value = 1.041855863584007
$ws.Range("N17").Value = 1.012515056898161
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.019212430090238
$ws.Range("D18").Value = 1.024062889688603
$ws.Range("E18").Value = 1.029590663943324
$ws.Range("F18").Value = 1.038750100963661
$ws.Range("I18").Value = 1.027846475420643
$ws.Range("J18").Value = 1.025428432187547
$ws.Range("K18").Value = 1.027428402284406
$ws.Range("L18").Value = 1.032936790721126
$ws.Range("M18").Value = 1.042064605412822
$ws.Range("N18").Value = 1.012554340693187
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.01928081587741
$ws.Range("D19").Value = 1.024111502742207
$ws.Range("E19").Value = 1.029656930106545
$ws.Range("F19").Value = 1.038835841358549
$ws.Range("I19").Value = 1.027856536367336
$ws.Range("J19").Value = 1.025468875289347
$ws.Range("K19").Value = 1.027462348560192
$ws.Range("L19").Value = 1.032988419576379
$ws.Range("M19").Value = 1.042135791210806
$ws.Range("N19").Value = 1.012567731909826
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.01897498974013
$ws.Range("D20").Value = 1.023894060467789
$ws.Range("E20").Value = 1.029360619535346
$ws.Range("F20").Value = 1.038452476883963
$ws.Range("I20").Value = 1.027811384119273
$ws.Range("J20").Value = 1.025287966616192
$ws.Range("K20").Value = 1.027310444502654
$ws.Range("L20").Value = 1.0327575195622
$ws.Range("M20").Value = 1.041817471955899
$ws.Range("N20").Value = 1.012507829260667
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.017980920182731
$ws.Range("D21").Value = 1.02318655734239
$ws.Range("E21").Value = 1.02839812065497
$ws.Range("F21").Value = 1.037207648441028
$ws.Range("I21").Value = 1.027661831025151
$ws.Range("J21").Value = 1.024699158805494
$ws.Range("K21").Value = 1.026815043948112
$ws.Range("L21").Value = 1.032006781606406
$ws.Range("M21").Value = 1.040783297968024
$ws.Range("N21").Value = 1.012312836003719
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.017356068392842
$ws.Range("D22").Value = 1.022741287726856
$ws.Range("E22").Value = 1.027793605927109
$ws.Range("F22").Value = 1.036426150005383
$ws.Range("I22").Value = 1.027565690475808
$ws.Range("J22").Value = 1.024328453884678
$ws.Range("K22").Value = 1.026502387084785
$ws.Range("L22").Value = 1.031534722010906
$ws.Range("M22").Value = 1.040133622637817
$ws.Range("N22").Value = 1.012190050782426
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.017687298717018
$ws.Range("D23").Value = 1.022977374216736
$ws.Range("E23").Value = 1.028114008941824
$ws.Range("F23").Value = 1.036840325480259
$ws.Range("I23").Value = 1.027616856198664
$ws.Range("J23").Value = 1.024525018471667
$ws.Range("K23").Value = 1.026668243761449
$ws.Range("L23").Value = 1.031784973288884
$ws.Range("M23").Value = 1.040477975304275
$ws.Range("N23").Value = 1.012255158993906
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.018991658071525
$ws.Range("D24").Value = 1.023905914372144
$ws.Range("E24").Value = 1.029376766828899
$ws.Range("F24").Value = 1.038473366446127
$ws.Range("I24").Value = 1.027813855560131
$ws.Range("J24").Value = 1.0252978295435
$ws.Range("K24").Value = 1.027318729888819
$ws.Range("L24").Value = 1.032770105021301
$ws.Range("M24").Value = 1.041834819293297
$ws.Range("N24").Value = 1.012511095183306
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.02050595312487
$ws.Range("D25").Value = 1.024981489057416
$ws.Range("E25").Value = 1.030844902080848
$ws.Range("F25").Value = 1.040373504768707
$ws.Range("I25").Value = 1.02803324576764
$ws.Range("J25").Value = 1.026192432829775
$ws.Range("K25").Value = 1.028068410558048
$ws.Range("L25").Value = 1.033913085691029
$ws.Range("M25").Value = 1.043411729855731
$ws.Range("N25").Value = 1.012807276805556
